# Update the jefaturas template:
#  - Replace the support email address shown (and linked) in column J
#    (correo_soporte) from "SeguridadTI@clarovtr.cl" to
#    "mariavyeguezp@gmail.com", and turn those cells into real mailto
#    hyperlinks (matching the existing "curso_link" hyperlink style).
#  - Update the sheet view's selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newEmail = "mariavyeguezp@gmail.com"
$mailto = "mailto:" + $newEmail

# 1) Update the cell values first so the old shared string
#    ("SeguridadTI@clarovtr.cl") becomes unused and the new one gets
#    appended to the shared strings table.
$ws.Range("J2").Value2 = $newEmail
$ws.Range("J3").Value2 = $newEmail
$ws.Range("J4").Value2 = $newEmail

# 2) Turn J2 into a mailto hyperlink (single cell). The cell text already
#    equals the target address so no explicit TextToDisplay is needed.
$ws.Hyperlinks.Add($ws.Range("J2"), $mailto)

# 3) Turn J3:J4 into a single merged mailto hyperlink reference, keeping
#    the visible text as the e-mail address.
$ws.Hyperlinks.Add($ws.Range("J3:J4"), $mailto, [Type]::Missing, [Type]::Missing, $newEmail)

# 4) Re-apply the existing "hyperlink" cell style (same style already used
#    by I2 / curso_link) to J2:J4 so they visually match and reuse the
#    same style index instead of creating new duplicate styles.
$hyperlinkStyle = $ws.Range("I2").Style
$ws.Range("J2:J4").Style = $hyperlinkStyle

# 5) Update the sheet view: scroll so column E is the left-most visible
#    column, and move the active selection to J12.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J12").Select() | Out-Null
